# Theme test case: add test cases for Normal theme
#
# The original sheet had an unused, completely empty row 1 above the
# header row. Remove it so the data starts at row 1 (this shifts every
# row up by one and Excel auto-adjusts the SUM() formulas), then fix a
# few "bouns" -> "bonus" typos in column I section headers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the blank leading row; everything below shifts up by one row.
$ws.Rows("1:1").Delete()

# Correct the "bouns" -> "bonus" typos in the section header notes.
# (These were at I3/I12/I22 before the delete; now at I2/I11/I21.)
$ws.Range("I2").Value = "no coin/diamond bonus"
$ws.Range("I11").Value = "with coin/diamond bonus"
$ws.Range("I21").Value = "coin/diamond oak + bonus"

# Leave the selection where the author ended up after the edit.
$ws.Range("I21").Select()
